$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1890
$ws1.Range("F3").Value = 502
$ws1.Range("F6").Value = 2647
$ws1.Range("F7").Value = 180
$ws1.Range("F9").Value = 178
$ws1.Range("F10").Value = 1557
$ws1.Range("F11").Value = 542
$ws1.Range("F13").Value = 338
$ws1.Range("F21").Value = 193
$ws1.Range("F23").Value = 1695
$ws1.Range("F25").Value = 414
$ws1.Range("F26").Value = 29

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1890
$ws4.Range("F4").Value = 502
$ws4.Range("F5").Value = 0
$ws4.Range("F7").Value = 2647
$ws4.Range("F8").Value = 180
$ws4.Range("F10").Value = 178
$ws4.Range("F11").Value = 1557
$ws4.Range("F12").Value = 542
$ws4.Range("F14").Value = 338
$ws4.Range("F22").Value = 193
$ws4.Range("F24").Value = 1695
$ws4.Range("F26").Value = 414
$ws4.Range("F27").Value = 29
